# Update price, config and deal data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (set 31142) lost its Collection value ("N/A" -> blank)
$ws.Range("D12").Value = ""

# A new LEGO set (31154) was added to the catalogue; it belongs right
# after row 12 (sorted by ID_Set), so insert a fresh row at position 13
# and push every row currently at 13..28 down to 14..29.
$ws.Rows.Item(13).Insert()

# Columns A (ID_Set) and C (nbPieces) are stored as text in this sheet,
# so force text formatting before writing numeric-looking values.
$ws.Range("A13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"

$ws.Range("A13").Value = "31154"
$ws.Range("B13").Value = "Les animaux de la forêt : Le renard roux"
$ws.Range("C13").Value = "667"
$ws.Range("D13").Value = "N/A"
$ws.Range("E13").Value = "https://www.lego.com/cdn/cs/set/assets/blt73662553f9401b9b/31154.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1"
$ws.Range("F13").Value = "https://www.lego.com/fr-fr/product/31154"
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""
